$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear existing hyperlinks so we can rebuild a clean, fully-sequential set
# (row positions are being reshuffled, so old hyperlink row refs would go stale)
$ws.Hyperlinks.Delete()

$rows = @(
    @{ A="2025-10-06 12:37:21"; B="【AI開発者募集】多機能転売ツールの構築をお願いします!"; C="システム開発"; D="500,000 円 ~ 1,000,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5407785"; G=420; H="🔥AI,Ai ◆ツール,開発" }
    @{ A="2025-10-06 12:37:21"; B="【生成AI】マーケティング事業の新プロダクト開発相談"; C="システム開発"; D="100,000 円 ~ 200,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5407690"; G=368; H="🔥AI,Ai ◆開発" }
    @{ A="2025-10-06 12:37:21"; B="【AI構築】Gmail助言テキストの自動記録システム構築依頼"; C="システム開発"; D="300,000 円 ~ 500,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5407751"; G=325; H="🔥AI,Ai" }
    @{ A="2025-10-06 12:37:21"; B="あなたAIクローン構築パートナー募集・モデル制作&新規依頼"; C="システム開発"; D="100,000 円 ~ 200,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5399534"; G=303; H="🔥AI,Ai" }
    @{ A="2025-10-06 12:37:21"; B="InstagramのAPIをどのように取得するかが不明なので代行して欲しい"; C="システム開発"; D="10,000 円 ~ 20,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5407720"; G=180; H="🔥API" }
    @{ A="2025-10-06 12:37:21"; B="初回 【急募】大規模データ収集自動化(スクレイピング・DB連携・エラー管理)案件"; C="システム開発"; D="50,000 円 ~ 100,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5407281"; G=158; H="◆自動化,スクレイピング ◇管理" }
    @{ A="2025-10-06 12:37:21"; B="【フルリモート】WordPress開発スタッフ募集"; C="システム開発"; D="50,000 円 ~ 100,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5407811"; G=88; H="◆開発 ○WordPress" }
    @{ A="2025-10-06 12:37:21"; B="見積書の計算から、社内受注表および受領証までを一括で作成できるシステム"; C="システム開発"; D="50,000 円 ~ 100,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5407766"; G=33; H=$null }
    @{ A="2025-10-06 12:37:21"; B="Access業務システムのクラウド化(ZOHO Creator使用)をお手伝いください!"; C="システム開発"; D="200,000 円 ~ 300,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5407746"; G=33; H=$null }
    @{ A="2025-10-06 12:37:21"; B="GCP上で動かしているWebサーバーの停止について"; C="システム開発"; D="~ 5,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5407840"; G=10; H=$null }
    @{ A="2025-10-06 12:37:21"; B="【急募】AWSマイクロサービスのデバッグ・最適化支援者募集"; C="システム開発"; D="5,000 円 ~ 10,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5407390"; G=10; H=$null }
    @{ A="2025-10-06 12:37:21"; B="【急募】サーバー移転後のWelcartクレカ決済不具合解消依頼"; C="システム開発"; D="8,000 円 ~ 9,000 円 / 募集期間 2 日、取引期間 0 日"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5407516"; G=10; H=$null }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $fcell = $ws.Cells.Item($r, 6)
    $fcell.Value = $row.F
    $ws.Hyperlinks.Add($fcell, $row.F)
    $fcell.Style = "Hyperlink"
    $ws.Cells.Item($r, 7).Value = $row.G
    if ($row.H -ne $null) {
        $ws.Cells.Item($r, 8).Value = $row.H
    } else {
        $ws.Cells.Item($r, 8).Value = ""
    }
    $r = $r + 1
}

# Column width adjustments (engine adds a constant +5/6 "padding" offset when
# round-tripping ColumnWidth -> stored OOXML width, so we pre-compensate).
$ws.Columns.Item(2).ColumnWidth = 46 - 0.8333333333333334
$ws.Columns.Item(4).ColumnWidth = 39 - 0.8333333333333334

